# "Updated quick quote for both qa and stg"
#
# The quick-quote lookup cells on the Input sheet (B2, B3, B5) get
# refreshed to point at a new set of generated product codes. Only the
# cell *values* change - the cell formatting (style index) is left
# exactly as it was, so B2/B3/B5 keep referencing the same visual style
# they already had.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B2").Value = "prodXvSh"
$ws.Range("B3").Value = "prodPcCp"
$ws.Range("B5").Value = "prodpLJJ"
